$wb = $excel.ActiveWorkbook

# 1. Add the new "cell-dive" option to the "assay_type list" sheet (A2),
#    alongside the existing "Cell DIVE" value in A1.
$listWs = $wb.Worksheets.Item("assay_type list")
$listWs.Range("A2").Value = "cell-dive"

# 2. Update the data validation on the "assay_type" column (L) of the
#    "Export as TSV" sheet so it references the expanded list range and
#    reports both allowed values in its error message.
$mainWs = $wb.Worksheets.Item("Export as TSV")
$dv = $mainWs.Range("L2:L1048576").Validation
$dv.Modify(3, 1, 1, "'assay_type list'!`$A`$1:`$A`$2")
$dv.ErrorMessage = "Value must be one of: Cell DIVE / cell-dive."
